$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new blank rows starting at row 495, pushing the existing
# rows 495-505 down to 501-511 (dimension grows from A1:T505 to A1:T511).
$ws.Rows("495:500").Insert()

# Common column values shared by every row in this block.
$mercadoId = 5
$mercado   = "Macroferia Regional de Talca"
$region    = "Maule"
$codreg    = 7
$tipo      = "Fruta"
$productoId = 100103
$producto  = "Frutos de hueso (carozo)"
$categoriaId = 100103006
$categoria = "Nectarín"
$origen    = "Región de O'Higgins"

$newRows = @(
  @{ Row=495; D=44595; K='Artic Star';  L='Especial'; M=300; N=12000; O=12000; P=12000; Q='$/bandeja 15 kilos granel'; S=800; T=15 },
  @{ Row=496; D=44595; K='Artic Star';  L='Primera';  M=300; N=9000;  O=9000;  P=9000;  Q='$/bandeja 15 kilos granel'; S=600; T=15 },
  @{ Row=497; D=44595; K='Red Diamond'; L='Especial'; M=300; N=13000; O=13000; P=13000; Q='$/bandeja 15 kilos granel'; S=867; T=15 },
  @{ Row=498; D=44595; K='Red Diamond'; L='Primera';  M=300; N=10000; O=10000; P=10000; Q='$/bandeja 15 kilos granel'; S=667; T=15 },
  @{ Row=499; D=44595; K='Sun Rise';    L='Especial'; M=300; N=14000; O=14000; P=14000; Q='$/bandeja 15 kilos granel'; S=933; T=15 },
  @{ Row=500; D=44595; K='Sun Rise';    L='Primera';  M=300; N=12000; O=12000; P=12000; Q='$/bandeja 15 kilos granel'; S=800; T=15 }
)

foreach ($r in $newRows) {
  $row = $r.Row
  $ws.Cells.Item($row, 1).Value  = $mercadoId
  $ws.Cells.Item($row, 2).Value  = $mercado
  $ws.Cells.Item($row, 3).Value  = $region
  $ws.Cells.Item($row, 4).Value  = $r.D
  $ws.Cells.Item($row, 5).Value  = $codreg
  $ws.Cells.Item($row, 6).Value  = $tipo
  $ws.Cells.Item($row, 7).Value  = $productoId
  $ws.Cells.Item($row, 8).Value  = $producto
  $ws.Cells.Item($row, 9).Value  = $categoriaId
  $ws.Cells.Item($row, 10).Value = $categoria
  $ws.Cells.Item($row, 11).Value = $r.K
  $ws.Cells.Item($row, 12).Value = $r.L
  $ws.Cells.Item($row, 13).Value = $r.M
  $ws.Cells.Item($row, 14).Value = $r.N
  $ws.Cells.Item($row, 15).Value = $r.O
  $ws.Cells.Item($row, 16).Value = $r.P
  $ws.Cells.Item($row, 17).Value = $r.Q
  $ws.Cells.Item($row, 18).Value = $origen
  $ws.Cells.Item($row, 19).Value = $r.S
  $ws.Cells.Item($row, 20).Value = $r.T
}
